$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Extend the leave table (Table1) from A8:K131 to A8:K137 by inserting six
#    new rows right before the old "last row" (row 131). Excel shifts the old
#    row 131 (with its distinct bottom-border style) down to row 137, and the
#    newly inserted rows 131-136 inherit generic default formatting that we
#    will fix up below by pasting the formats from row 130 (a normal body row).
# ---------------------------------------------------------------------------
$ws.Range("A131:A136").EntireRow.Insert()

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K137"))

$ws.Range("A130:K130").Copy()
$ws.Range("A131:K136").PasteSpecial(-4122)

# Restore the per-cell "EARNED " mirror formula (written individually so each
# row keeps its own literal formula instead of Excel collapsing them into one
# shared-formula group).
for ($r = 131; $r -le 137; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
}

# ---------------------------------------------------------------------------
# 2) New leave entries recorded in rows 100-107.
# ---------------------------------------------------------------------------
$ws.Range("B100").Value = "SL(2-0-00)"
$ws.Cells.Item(100, 8).Value = 2
$ws.Range("K100").Value = "12/28,29/2022"
$ws.Range("K99").Copy()
$ws.Range("K100").PasteSpecial(-4122)

$ws.Range("A101").NumberFormat = "@"
$ws.Range("A101").Value = "2023"
$ws.Range("A100").Copy()
$ws.Range("A101").PasteSpecial(-4122)
$ws.Range("K99").Copy()
$ws.Range("K101").PasteSpecial(-4122)

$ws.Range("A102").Value = 44927
$ws.Range("B102").Value = "SL(1-0-00)"
$ws.Range("C102").Value = 1.25
$ws.Cells.Item(102, 8).Value = 1
$ws.Range("K102").Value = 44929
$ws.Range("K99").Copy()
$ws.Range("K102").PasteSpecial(-4122)

$ws.Range("B103").Value = "SL(3-0-0)"
$ws.Cells.Item(103, 8).Value = 3
$ws.Range("K103").Value = "1/18,19,20/2023"
$ws.Range("K99").Copy()
$ws.Range("K103").PasteSpecial(-4122)

$ws.Range("A104").Value = 44958
$ws.Range("B104").Value = "VL(4-0-0)"
$ws.Range("C104").Value = 1.25
$ws.Range("D104").Value = 4
$ws.Range("K104").Value = "2/10-15/2023"

$ws.Range("B105").Value = "SP(3-0-00)"
$ws.Range("K105").Value = "2/7-9/2023"

$ws.Range("B106").Value = "VL(5-0-0)"
$ws.Range("D106").Value = 5
$ws.Range("K106").Value = "1/30,31, 2/1,2,3"

$ws.Range("B107").Value = "SL(1-0-00)"
$ws.Cells.Item(107, 8).Value = 1
$ws.Range("K107").Value = 44952
$ws.Range("K99").Copy()
$ws.Range("K107").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) PERIOD (column A) end-of-month dates for the remaining periods, rows
#    108-133 (the table grew enough to cover leave periods through 2025).
# ---------------------------------------------------------------------------
$periodDates = @{
    108 = 44986; 109 = 45017; 110 = 45047; 111 = 45078; 112 = 45108
    113 = 45139; 114 = 45170; 115 = 45200; 116 = 45231; 117 = 45261
    118 = 45292; 119 = 45323; 120 = 45352; 121 = 45383; 122 = 45413
    123 = 45444; 124 = 45474; 125 = 45505; 126 = 45536; 127 = 45566
    128 = 45597; 129 = 45627; 130 = 45658; 131 = 45689; 132 = 45717
    133 = 45748
}
foreach ($r in $periodDates.Keys) {
    $ws.Cells.Item($r, 1).Value = $periodDates[$r]
}

# ---------------------------------------------------------------------------
# 4) View state: move the active selection to B108 (matches the author's
#    last position after recording the new leave periods).
# ---------------------------------------------------------------------------
$ws.Range("B108").Select()
